$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.538.67'
$ws.Range("E2").Value = '  +0.07%  '

$ws.Range("D3").Value = '1.815.37'
$ws.Range("E3").Value = '  +0.13%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("E5").Value = '  -1.00%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.595'
$ws.Range("E6").Value = '  +2.51%  '

$ws.Range("E7").Value = '  -0.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '38.34'
$ws.Range("E8").Value = '  +5.92%  '

$ws.Range("E9").Value = '  -4.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0680'
$ws.Range("E10").Value = '  -3.03%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0977'
$ws.Range("E11").Value = '  +1.40%  '

$ws.Range("D12").Value = '2.076.48'
$ws.Range("E12").Value = '  -0.01%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.32'
$ws.Range("E13").Value = '  -1.72%  '

$ws.Range("D14").Value = '1.819.27'
$ws.Range("E14").Value = '  +0.37%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.637'
$ws.Range("E15").Value = '  -2.40%  '

$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '34.508.00'
$ws.Range("E16").Value = '  -0.02%  '

$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.45'
$ws.Range("E17").Value = '  -2.38%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.73'
$ws.Range("E18").Value = '  -1.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.90'
$ws.Range("E19").Value = '  -1.91%  '

$ws.Range("E20").Value = '  -3.31%  '

$ws.Range("E21").Value = '  -2.94%  '

$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.16'
$ws.Range("E23").Value = '  -2.00%  '

$ws.Range("E24").Value = '  +4.18%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.06'
$ws.Range("E25").Value = '  -0.86%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.94'
$ws.Range("E26").Value = '  -1.97%  '

$ws.Range("E27").Value = '  +4.74%  '

$ws.Range("E28").Value = '  +0.58%  '

$ws.Range("E29").Value = '  -0.07%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.81'
$ws.Range("E30").Value = '  -2.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.23'
$ws.Range("E31").Value = '  -1.40%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0521'
$ws.Range("E32").Value = '  -2.69%  '

$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.88'
$ws.Range("E33").Value = '  -5.32%  '

$ws.Range("E34").Value = '  -1.00%  '

$ws.Range("D35").Value = '1.365.49'
$ws.Range("E35").Value = '  -2.80%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.650'
$ws.Range("E36").Value = '  -5.38%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.06'
$ws.Range("E37").Value = '  -1.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.36'
$ws.Range("E38").Value = '  -7.48%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0188'
$ws.Range("E39").Value = '  -2.62%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.958'
$ws.Range("E40").Value = '  -1.74%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.44'
$ws.Range("E41").Value = '  +1.37%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '81.97'
$ws.Range("E42").Value = '  -3.36%  '

$ws.Range("E43").Value = '  +0.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.78'
$ws.Range("E44").Value = '  -1.75%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.89'
$ws.Range("E45").Value = '  +2.96%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0509'
$ws.Range("E46").Value = '  +0.53%  '

$ws.Range("D47").Value = '1.977.20'
$ws.Range("E47").Value = '  +0.05%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.81'
$ws.Range("E48").Value = '  -4.45%  '

$ws.Range("E49").Value = '  -0.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '102.33'
$ws.Range("E50").Value = '  -3.69%  '

$ws.Range("D51").Value = '0.0₆0123'
$ws.Range("E51").Value = '  -5.18%  '
